$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.794.01"
$ws.Range("E2").Value = "  +1.87%  "
$ws.Range("D3").Value = "3.275.88"
$ws.Range("E3").Value = "  -0.81%  "
$ws.Range("D4").Value = "'0.997"
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").Value = "'569.41"
$ws.Range("D6").Value = "'175.96"
$ws.Range("E6").Value = "  -4.11%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.580"
$ws.Range("E8").Value = "  +1.70%  "
$ws.Range("D9").Value = "3.271.68"
$ws.Range("E9").Value = "  -0.77%  "
$ws.Range("E10").Value = "  -1.63%  "
$ws.Range("D11").Value = "'0.571"
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("D12").Value = "'45.69"
$ws.Range("E12").Value = "  -2.15%  "
$ws.Range("E13").Value = "  +1.72%  "
$ws.Range("D14").Value = "'700.72"
$ws.Range("E14").Value = "  +10.20%  "
$ws.Range("D15").Value = "3.799.18"
$ws.Range("E15").Value = "  -0.75%  "
$ws.Range("D16").Value = "'8.31"
$ws.Range("E16").Value = "  -1.68%  "
$ws.Range("D17").Value = "66.869.70"
$ws.Range("E17").Value = "  +1.72%  "
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("D19").Value = "3.272.87"
$ws.Range("E19").Value = "  -0.77%  "
$ws.Range("E20").Value = "  -2.58%  "
$ws.Range("D21").Value = "'10.73"
$ws.Range("E21").Value = "  -2.47%  "
$ws.Range("D22").Value = "'0.888"
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("D23").Value = "'16.91"
$ws.Range("E23").Value = "  -5.56%  "
$ws.Range("E24").Value = "  +3.32%  "
$ws.Range("D25").Value = "'97.97"
$ws.Range("E25").Value = "  -3.00%  "
$ws.Range("D26").Value = "'3.87"
$ws.Range("E26").Value = "  -2.30%  "
$ws.Range("E27").Value = "  -1.34%  "
$ws.Range("D28").Value = "'9.32"
$ws.Range("E28").Value = "  -0.92%  "
$ws.Range("D29").Value = "'32.86"
$ws.Range("E29").Value = "  +6.32%  "
$ws.Range("D30").Value = "'8.43"
$ws.Range("E30").Value = "  +0.73%  "
$ws.Range("D31").Value = "'6.81"
$ws.Range("E31").Value = "  +4.31%  "
$ws.Range("D32").Value = "'566.64"
$ws.Range("E32").Value = "  -3.05%  "
$ws.Range("D33").Value = "3.886.76"
$ws.Range("E33").Value = "  +1.07%  "
$ws.Range("D34").Value = "'10.81"
$ws.Range("E34").Value = "  -0.70%  "
$ws.Range("E35").Value = "  -1.22%  "
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").Value = "'55.49"
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").Value = "'3.33"
$ws.Range("E38").Value = "  -9.89%  "
$ws.Range("E39").Value = "  +1.48%  "
$ws.Range("D40").Value = "'2.61"
$ws.Range("E40").Value = "  +0.48%  "
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "'31.84"
$ws.Range("E41").Value = "  -2.21%  "
$ws.Range("B42").Value = "ApeXProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D42").Value = "'3.34"
$ws.Range("E42").Value = "  -1.85%  "
$ws.Range("D43").Value = "0.0₃0674"
$ws.Range("E43").Value = "  -3.22%  "
$ws.Range("D44").Value = "'3.01"
$ws.Range("E44").Value = "  -3.39%  "
$ws.Range("D45").Value = "'0.328"
$ws.Range("E45").Value = "  -1.36%  "
$ws.Range("D46").Value = "'0.0405"
$ws.Range("E46").Value = "  -0.32%  "
$ws.Range("E47").Value = "  +0.60%  "
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("E49").Value = "  +0.53%  "
$ws.Range("E50").Value = "  +7.56%  "
$ws.Range("D51").Value = "'129.79"
$ws.Range("E51").Value = "  +0.06%  "
